$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.316.70"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.573.71"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'588.46"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").Value = "'186.42"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "3.562.82"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("D11").Value = "'0.646"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "'54.55"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "4.140.47"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'19.47"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "70.330.26"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "3.555.31"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "12.49"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'544.09"
$ws.Range("E21").Value = "  +10.55%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "17.98"
$ws.Range("E23").Value = "  -7.72%  "
$ws.Range("D24").Value = "'4.69"
$ws.Range("E24").Value = "  +7.95%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "'96.00"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "11.25"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'3.00"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "'9.14"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "'32.27"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "'7.35"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "'12.55"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").Value = "65.25"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "'555.98"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("D36").Value = "'3.22"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Value = "'0.416"
$ws.Range("E37").Value = "  +5.36%  "
$ws.Range("D38").Value = "'38.62"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "3.365.39"
$ws.Range("E43").Value = "  +4.18%  "
$ws.Range("E44").Value = "  -7.59%  "
$ws.Range("D45").Value = "'3.55"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.136"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'9.16"
$ws.Range("E49").Value = "  -5.48%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'137.42"
$ws.Range("E51").Value = "  +1.23%  "
